$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps explicit Text format so that
# numeric-looking strings (e.g. "241.79") are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.019.50"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.910.78"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "0.7847"
$ws.Range("E5").Value = "  +5.29%  "
$ws.Range("D6").Value = "241.79"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").Value = "26.10"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "0.06900"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").Value = "0.07957"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "1.906.40"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "0.7420"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "5.210"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "92.99"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "30.032.61"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "13.97"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "5.881"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("D19").Value = "246.18"
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").Value = "0.000007749"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "2.150.01"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "6.875"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("D25").Value = "169.08"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("D26").Value = "9.287"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "0.1376"
$ws.Range("E27").Value = "  +8.61%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "1.377"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "4.322"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").Value = "0.05475"
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").Value = "1.256"
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("D36").Value = "0.7342"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").Value = "0.01934"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").Value = "2.793"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("D41").Value = "0.4415"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").Value = "72.05"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "0.8374"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").Value = "1.877"
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("D46").Value = "100.40"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.522"
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.771"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "981.20"
$ws.Range("E49").Value = "  +8.51%  "
$ws.Range("D50").Value = "2.058.34"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "36.23"
$ws.Range("E51").Value = "  -1.20%  "
